# "Added form matrix mapping."
#
# Update the GA (row 2) and AZ (row 6) dynamic-fields test cases with their
# real form numbers / descriptions, tweak a couple of step/result strings,
# and tidy up the header row styling + column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2 : GA (Georgia) Manhole Liability Coverage test case ----
$ws.Range("B2").Value = "TC_1_Verify the dynamic fields for BAS UM (GA) 12 49 - Georgia Manhole Liability Coverage"
$ws.Range("C2").Value = "Verify the dynamic fields for BAS UM (GA) 12 49 - Georgia Manhole Liability Coverage"
$ws.Range("E2").Value = "Login to PC and initiate a submission for GA"

# ---- Row 4 & Row 8 : coverage-term step (shared text, both rows updated) ----
$ws.Range("E4").Value = "Add  will triggered when Manhole Liability coverage is selected"
$ws.Range("F4").Value = "User should be able to add all the terms"
$ws.Range("E8").Value = "Add  will triggered when Manhole Liability coverage is selected"
$ws.Range("F8").Value = "User should be able to add all the terms"

# ---- Row 6 : AZ (Arizona) Manhole Liability Coverage test case ----
$ws.Range("B6").Value = "TC_1_Verify the dynamic fields for BAS UM (AZ) 03 12 - Arizona Manhole Liability Coverage"
$ws.Range("C6").Value = "Verify the dynamic fields for BAS UM (AZ) 03 12 - Arizona Manhole Liability Coverage"
$ws.Range("E6").Value = "Login to PC and initiate a submission for AZ"

# ---- Header row: drop the (colourless) solid-fill shading, keep Bold ----
$ws.Range("A1:G1").Interior.Pattern = -4142

# ---- Widen Summary/Description columns to match ----
$ws.Columns.Item(2).ColumnWidth = 49.17
$ws.Columns.Item(3).ColumnWidth = 49.17
